$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 36000
$ws.Range("J75").Value = 36000
$ws.Range("L75").Value = 36000
$ws.Range("N75").Value = -37872

$ws.Range("H76").Value = 2648688.5
$ws.Range("I76").Value = 3089603
$ws.Range("J76").Value = 3202
$ws.Range("K76").Value = 3089603
$ws.Range("L76").Value = 3202
$ws.Range("M76").Value = -3089288
$ws.Range("N76").Value = -3832

$ws.Range("H78").Value = 36000
$ws.Range("J78").Value = 36000
$ws.Range("L78").Value = 108000
$ws.Range("N78").Value = -117360

$ws.Range("H79").Value = 2648688.5
$ws.Range("I79").Value = 3089603
$ws.Range("J79").Value = 3202
$ws.Range("K79").Value = 3089603
$ws.Range("L79").Value = 3202
$ws.Range("M79").Value = -3088511
$ws.Range("N79").Value = -5386

$ws.Range("H135").Value = 4159.15
$ws.Range("I135").Value = 2639.0715
$ws.Range("J135").Value = 7706
$ws.Range("K135").Value = 23751.6435
$ws.Range("L135").Value = 69354
$ws.Range("M135").Value = -21216.6435
$ws.Range("N135").Value = -74424

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 1361.7428
$ws.Range("I137").Value = 872.7273
$ws.Range("J137").Value = 2189.3076
$ws.Range("K137").Value = 2618.1819
$ws.Range("L137").Value = 6567.9228
$ws.Range("M137").Value = -68.18190000000004
$ws.Range("N137").Value = -11667.9228

$ws.Range("H139").Value = 77840
$ws.Range("J139").Value = 77840
$ws.Range("L139").Value = 77840
$ws.Range("N139").Value = -88120

$ws.Range("H140").Value = 70560.87
$ws.Range("J140").Value = 96635.71000000001
$ws.Range("L140").Value = 96635.71000000001
$ws.Range("N140").Value = -106995.71

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 748436.5600000001
$ws.Range("I32").Value = 9834.816000000001
$ws.Range("J32").Value = 4367585
$ws.Range("K32").Value = 9834.816000000001
$ws.Range("L32").Value = 4367585
$ws.Range("M32").Value = -9547.816000000001
$ws.Range("N32").Value = -4368159

$ws.Range("H61").Value = 8420.666999999999
$ws.Range("I61").Value = 9282
$ws.Range("J61").Value = 4114
$ws.Range("K61").Value = 9282
$ws.Range("L61").Value = 4114
$ws.Range("M61").Value = -9070
$ws.Range("N61").Value = -4538

$ws.Range("H74").Value = 862.4
$ws.Range("I74").Value = 862.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 862.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 11.60000000000002
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 862.4
$ws.Range("I77").Value = 862.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4312
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 56
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 5412.276
$ws.Range("I122").Value = 5789.5713
$ws.Range("J122").Value = 4421.875
$ws.Range("K122").Value = 17368.7139
$ws.Range("L122").Value = 13265.625
$ws.Range("M122").Value = -14918.7139
$ws.Range("N122").Value = -18165.625

$ws.Range("H136").Value = 8420.666999999999
$ws.Range("I136").Value = 9282
$ws.Range("J136").Value = 4114
$ws.Range("K136").Value = 27846
$ws.Range("L136").Value = 12342
$ws.Range("M136").Value = -25296
$ws.Range("N136").Value = -17442

$ws.Range("H140").Value = 103032.9
$ws.Range("J140").Value = 103032.9
$ws.Range("L140").Value = 103032.9
$ws.Range("N140").Value = -113392.9

$ws.Range("H141").Value = 65330
$ws.Range("J141").Value = 65330
$ws.Range("L141").Value = 65330
$ws.Range("N141").Value = -75690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2580
$ws.Range("I99").Value = 1367.25
$ws.Range("K99").Value = 1367.25
$ws.Range("M99").Value = 130.75

$ws.Range("H128").Value = 1750
$ws.Range("I128").Value = 1750
$ws.Range("K128").Value = 5250
$ws.Range("M128").Value = -2760

$ws.Range("H134").Value = 1271
$ws.Range("I134").Value = 1012
$ws.Range("J134").Value = 1322.8
$ws.Range("K134").Value = 3036
$ws.Range("L134").Value = 3968.4
$ws.Range("M134").Value = -501
$ws.Range("N134").Value = -9038.4

$ws.Range("H138").Value = 69833.336
$ws.Range("J138").Value = 69833.336
$ws.Range("L138").Value = 69833.336
$ws.Range("N138").Value = -80113.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 28572640
$ws.Range("I122").Value = 1496.4
$ws.Range("J122").Value = 100000500
$ws.Range("K122").Value = 4489.200000000001
$ws.Range("L122").Value = 300001500
$ws.Range("M122").Value = -2039.200000000001
$ws.Range("N122").Value = -300006400

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws.Range("H140").Value = 89900
$ws.Range("J140").Value = 89900
$ws.Range("L140").Value = 89900
$ws.Range("N140").Value = -100260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 458.33334
$ws.Range("I122").Value = 356.48386
$ws.Range("J122").Value = 1089.8
$ws.Range("K122").Value = 3208.35474
$ws.Range("L122").Value = 9808.199999999999
$ws.Range("M122").Value = -758.3547399999998
$ws.Range("N122").Value = -14708.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 3831.6667
$ws.Range("I122").Value = 2897.4
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 8692.200000000001
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -6242.200000000001
$ws.Range("N122").Value = -19898.5

$ws.Range("H132").Value = 3561.2666
$ws.Range("I132").Value = 2545.3333
$ws.Range("J132").Value = 5085.1665
$ws.Range("K132").Value = 7635.999899999999
$ws.Range("L132").Value = 15255.4995
$ws.Range("M132").Value = -5105.999899999999
$ws.Range("N132").Value = -20315.4995

$ws.Range("H138").Value = 69300
$ws.Range("J138").Value = 69300
$ws.Range("L138").Value = 69300
$ws.Range("N138").Value = -79580

$ws.Range("H139").Value = 52431.5
$ws.Range("J139").Value = 52431.5
$ws.Range("L139").Value = 52431.5
$ws.Range("N139").Value = -62711.5

$ws.Range("H140").Value = 89784.5
$ws.Range("J140").Value = 89784.5
$ws.Range("L140").Value = 89784.5
$ws.Range("N140").Value = -100144.5

$ws.Range("H141").Value = 47993.332
$ws.Range("J141").Value = 47993.332
$ws.Range("L141").Value = 47993.332
$ws.Range("N141").Value = -58353.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2049.7932
$ws.Range("I68").Value = 1972
$ws.Range("J68").Value = 2222.6667
$ws.Range("K68").Value = 1972
$ws.Range("L68").Value = 2222.6667
$ws.Range("M68").Value = -1223
$ws.Range("N68").Value = -3720.6667

$ws.Range("H71").Value = 2049.7932
$ws.Range("I71").Value = 1972
$ws.Range("J71").Value = 2222.6667
$ws.Range("K71").Value = 9860
$ws.Range("L71").Value = 11113.3335
$ws.Range("M71").Value = -6116
$ws.Range("N71").Value = -18601.3335

$ws.Range("H136").Value = 5158.5454
$ws.Range("I136").Value = 3513.3076
$ws.Range("J136").Value = 7535
$ws.Range("K136").Value = 10539.9228
$ws.Range("L136").Value = 22605
$ws.Range("M136").Value = -7989.9228
$ws.Range("N136").Value = -27705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1237.3846
$ws.Range("I136").Value = 963.65
$ws.Range("J136").Value = 1525.5264
$ws.Range("K136").Value = 2890.95
$ws.Range("L136").Value = 4576.5792
$ws.Range("M136").Value = -340.9499999999998
$ws.Range("N136").Value = -9676.5792

$ws.Range("H138").Value = 76740
$ws.Range("J138").Value = 87175
$ws.Range("L138").Value = 87175
$ws.Range("N138").Value = -97455

$ws.Range("H139").Value = 53376.875
$ws.Range("J139").Value = 53376.875
$ws.Range("L139").Value = 53376.875
$ws.Range("N139").Value = -63656.875

$ws.Range("H140").Value = 59800
$ws.Range("J140").Value = 59800
$ws.Range("L140").Value = 59800
$ws.Range("N140").Value = -70160

$ws.Range("H141").Value = 78653.44500000001
$ws.Range("J141").Value = 78653.44500000001
$ws.Range("L141").Value = 78653.44500000001
$ws.Range("N141").Value = -89013.44500000001
